$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 751 (2026/12/29 火 13 88 ...),
# shifting it and all following rows down by one. This grows the used
# range from A1:D792 to A1:D793.
$ws.Rows.Item(751).Insert()

# Populate the newly inserted row with the new daily entry.
# Column A holds a date-like label that must stay plain text (matching
# the rest of the sheet), so write it with a leading apostrophe and then
# clear the formatting Excel applies for the forced-text quote prefix.
$ws.Range("A751").Value = "'2026/01/30"
$ws.Range("A751").ClearFormats()
$ws.Range("B751").Value = "金"
$ws.Range("C751").Value = 20
$ws.Range("D751").Value = 201
